# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.1 = 49044.99 pesos`n✅ 49044.99 pesos = 12.07 = 977.52 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 82.65900000000001
$wsTasas.Range("O10").Value = 4054.01
$wsTasas.Range("N12").Value = 4063.99
$wsTasas.Range("O12").Value = 81
